$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A254:B254 already carry style "2" (the Book Antiqua / "常规 2" cell style used for the
# year/month columns of the data table). We use it as a format template below for the
# brand-new cells so they pick up the exact same style index instead of a freshly
# allocated one.

# --- Row 254 (already existed; keep its styles, just fill in the values) ---
$ws.Range("A254").Value = 2022
$ws.Range("B254").Value = 1
$ws.Range("C254").Value = 255.1369

# --- Row 255 (already existed; keep its styles, just fill in the values) ---
$ws.Range("A255").Value = 2022
$ws.Range("B255").Value = 2
$ws.Range("C255").Value = 299.0845

# --- Row 256 is brand new. Give A/B the "2" style; C keeps the default column style. ---
$ws.Range("A254:B254").Copy()
$ws.Range("A256:B256").PasteSpecial(-4122)
$ws.Range("A256").Value = 2022
$ws.Range("B256").Value = 3
$ws.Range("C256").Value = 305.6411

# --- Row 257 already existed (only A257, empty, style "4"); it keeps its
# ht="15" customHeight="1". Re-format A/B with style "2"; C gets the default style. ---
$ws.Range("A254:B254").Copy()
$ws.Range("A257:B257").PasteSpecial(-4122)
$ws.Range("A257").Value = 2022
$ws.Range("B257").Value = 4
$ws.Range("C257").Value = 401.5362

# --- Row 258 already existed, holding the footnote text in A258 (style "4").
# Stash that exact format in row 271 first so the footnote can be moved down there
# once this row is overwritten with numeric data. ---
$ws.Range("A258:A258").Copy()
$ws.Range("A271").PasteSpecial(-4122)

# Row 258 keeps its ht="15" customHeight="1". Re-format A/B with style "2"; C default.
$ws.Range("A254:B254").Copy()
$ws.Range("A258:B258").PasteSpecial(-4122)
$ws.Range("A258").Value = 2022
$ws.Range("B258").Value = 5
$ws.Range("C258").Value = 344.9735

# --- Row 259 already existed as a bare empty row (ht="15" customHeight="1", no
# cells). Give A/B style "2"; C default. ---
$ws.Range("A254:B254").Copy()
$ws.Range("A259:B259").PasteSpecial(-4122)
$ws.Range("A259").Value = 2022
$ws.Range("B259").Value = 6
$ws.Range("C259").Value = 310.6081

# --- Row 260 is brand new (no prior height). A/B style "2"; C default. ---
$ws.Range("A254:B254").Copy()
$ws.Range("A260:B260").PasteSpecial(-4122)
$ws.Range("A260").Value = 2022
$ws.Range("B260").Value = 7
$ws.Range("C260").Value = 237.5722

# --- Row 261 is brand new (no prior height). A/B style "2"; C default. ---
$ws.Range("A254:B254").Copy()
$ws.Range("A261:B261").PasteSpecial(-4122)
$ws.Range("A261").Value = 2022
$ws.Range("B261").Value = 8
$ws.Range("C261").Value = 384.6608

# --- Row 271: the footnote, moved down from its old home at A258. The format was
# already stashed there above (style "4", the font used for the citation note). ---
$ws.Range("A271").Value = "如需使用本指数进行相关研究，请注明原文来源：陈英楠、莫东翠、唐思华、李慧慧，《测量中国房地产政策不确定性研究》，《经济学》（季刊），2022年第22卷第2期，第405-424页。"
$ws.Rows(271).RowHeight = 14.65

$excel.CutCopyMode = 0

# --- Sheet view: drop the old scroll/selection state and select H8 instead ---
[void]$ws.Range("H8").Select()
